$wb = $excel.ActiveWorkbook

# xlPasteFormats = -4122 (used below to copy the bold/bordered header style
# from the last existing header cell onto the newly added header cell,
# reusing the workbook's existing style entry rather than creating a new one)
$xlPasteFormats = -4122

# --- Sheet: 对公产品台账202404 (add column N "当月新增户数") ---
$ws4 = $wb.Worksheets.Item("对公产品台账202404")

$ws4.Range("M1").Copy()
$ws4.Range("N1").PasteSpecial($xlPasteFormats)
$ws4.Range("N1").Value = "当月新增户数"

$n4values = @(0, 5, 0, 0, 0, 11, 0, 0, 0, 0, 0, 0, 0, 1, 17)
for ($i = 0; $i -lt $n4values.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 14).Value = $n4values[$i]
}

# --- Sheet: 个人经营贷202404 (add column P "当月新开户") ---
$ws5 = $wb.Worksheets.Item("个人经营贷202404")

$ws5.Range("O1").Copy()
$ws5.Range("P1").PasteSpecial($xlPasteFormats)
$ws5.Range("P1").Value = "当月新开户"

$p5values = @(0, 102, 38, 0, 37, 4, 5, 258, 261, 0, 3, 6, 0, 0, 714)
for ($i = 0; $i -lt $p5values.Length; $i++) {
    $row = $i + 2
    $ws5.Cells.Item($row, 16).Value = $p5values[$i]
}
